# salary-feb-2016.xlsx: add a new "Manager" sheet, give it a header row
# (Staff ID / Name, bold) and a staff id, and add one staff id row to each
# of the existing sheets (HR, Accounting, QA).

$wb = $excel.ActiveWorkbook

$hr = $wb.Worksheets.Item("HR")
$accounting = $wb.Worksheets.Item("Accounting")
$qa = $wb.Worksheets.Item("QA")

# --- fill in staff ids on the existing sheets -----------------------------
$hr.Range("A2").Value = 1
$hr.Range("A3").Value = 5

$accounting.Range("A2").Value = 2
# leave a stray selection on the Accounting sheet, as happened during the
# original editing session
[void]$accounting.Range("A20").Select()

$qa.Range("A2").Value = 3
[void]$qa.Range("B19").Select()

# --- create the new "Manager" sheet ---------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$manager = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$manager.Name = "Manager"

# reuse the same shared-string text already used for the other headers
$manager.Range("A1").Value = $hr.Range("A1").Text
$manager.Range("B1").Value = $hr.Range("B1").Text
$manager.Range("A1:B1").Font.Bold = $true

$manager.Range("A2").Value = 4

# --- leave the HR sheet/cell A2 selected & active, as in the final file ---
[void]$hr.Activate()
[void]$hr.Range("A2").Select()
